$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "remark" column (K),
# pushing remark from K -> M. This mirrors the commit that added
# "accountId" and "userId" columns to the transactions export.
$ws.Columns.Item(11).Insert()
$ws.Columns.Item(11).Insert()

# New header cells
$ws.Range("K1").Value2 = "accountId"
$ws.Range("L1").Value2 = "userId"

# New data cells per row
$ws.Range("K2").Value2 = "acc1"
$ws.Range("L2").Value2 = "AureliaSKY"

$ws.Range("K3").Value2 = "acc2"
$ws.Range("L3").Value2 = "AureliaSKY"

$ws.Range("K4").Value2 = "acc1"
$ws.Range("L4").Value2 = "AureliaSKY"

$ws.Range("K5").Value2 = "acc2"
$ws.Range("L5").Value2 = "AureliaSKY"

# Best-effort column widths to match the bestFit sizing from the source
# workbook (COM's ColumnWidth setter snaps to 1/6-character increments,
# so these are the closest achievable values).
$ws.Columns.Item(11).ColumnWidth = 8.833333333333334
$ws.Columns.Item(12).ColumnWidth = 9.833333333333334
